$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset")

# Rename antibodies from "Acme mAb N" to "COVIC N" for rows 2-11
for ($i = 1; $i -le 10; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "COVIC $i"
}

# Clear the qualitative measure results for rows 8-11 (Acme mAb 7-10 / COVIC 7-10)
$ws.Range("B8:B11").ClearContents()
